$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2-21 down to 3-22
$ws.Rows("2:2").Insert()
$ws.Range("A2:C2").ClearFormats()

# Populate the newly inserted row 2
$ws.Cells.Item(2, 1).Value = -0.0050396383740007
$ws.Cells.Item(2, 2).Value = -0.0062613687478005
$ws.Cells.Item(2, 3).Value = -0.0682641938328743

# Append 9 new data rows (23-31) at the bottom
$ws.Cells.Item(23, 1).Value = -0.4335615932941437
$ws.Cells.Item(23, 2).Value = 0.1406517177820205
$ws.Cells.Item(23, 3).Value = -0.8185594081878662
$ws.Cells.Item(24, 1).Value = 0.09498954564332961
$ws.Cells.Item(24, 2).Value = -0.7519751191139221
$ws.Cells.Item(24, 3).Value = -0.1093448773026466
$ws.Cells.Item(25, 1).Value = 0.1846340149641037
$ws.Cells.Item(25, 2).Value = -1.312596678733826
$ws.Cells.Item(25, 3).Value = 0.0687223374843597
$ws.Cells.Item(26, 1).Value = 0.6478226184844971
$ws.Cells.Item(26, 2).Value = -0.9091202020645142
$ws.Cells.Item(26, 3).Value = -0.1838704347610473
$ws.Cells.Item(27, 1).Value = -0.1064432710409164
$ws.Cells.Item(27, 2).Value = -0.09178250283002851
$ws.Cells.Item(27, 3).Value = 0.0652098655700683
$ws.Cells.Item(28, 1).Value = -0.042302418500185
$ws.Cells.Item(28, 2).Value = 0.3572034537792206
$ws.Cells.Item(28, 3).Value = 0.1937969923019409
$ws.Cells.Item(29, 1).Value = -0.2768746614456177
$ws.Cells.Item(29, 2).Value = 0.2338086664676666
$ws.Cells.Item(29, 3).Value = -0.1817324161529541
$ws.Cells.Item(30, 1).Value = 0.0734565481543541
$ws.Cells.Item(30, 2).Value = 0.1968513280153274
$ws.Cells.Item(30, 3).Value = 0.1055269688367843
$ws.Cells.Item(31, 1).Value = -0.1507309973239898
$ws.Cells.Item(31, 2).Value = -0.0175623763352632
$ws.Cells.Item(31, 3).Value = 0.08170322328805921
